# Auto-generated edit script applying the Marilith_Profits leve-profit recalculation update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 70
$ws.Range("H70").Value = 3361
$ws.Range("I70").Value = 2151
$ws.Range("J70").Value = 3547.1538
$ws.Range("K70").Value = 6453
$ws.Range("L70").Value = 10641.4614
$ws.Range("M70").Value = -6183
$ws.Range("N70").Value = -11181.4614

# ALC row 73
$ws.Range("H73").Value = 3361
$ws.Range("I73").Value = 2151
$ws.Range("J73").Value = 3547.1538
$ws.Range("K73").Value = 6453
$ws.Range("L73").Value = 10641.4614
$ws.Range("M73").Value = -5517
$ws.Range("N73").Value = -12513.4614

# ALC row 80
$ws.Range("H80").Value = 1799.5
$ws.Range("I80").Value = 600
$ws.Range("K80").Value = 1800
$ws.Range("M80").Value = -802

# ALC row 83
$ws.Range("H83").Value = 1799.5
$ws.Range("I83").Value = 600
$ws.Range("K83").Value = 5400
$ws.Range("M83").Value = -408

# ALC row 105
$ws.Range("H105").Value = 22000
$ws.Range("J105").Value = 22000
$ws.Range("L105").Value = 22000
$ws.Range("N105").Value = -28988

# ALC row 129
$ws.Range("H129").Value = 1627
$ws.Range("I129").Value = 833.5714
$ws.Range("K129").Value = 2500.7142
$ws.Range("M129").Value = 2499.2858

$ws = $wb.Worksheets.Item("ARM")
# ARM row 11
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

# ARM row 32
$ws.Range("H32").Value = 2531.691
$ws.Range("I32").Value = 2344.2075
$ws.Range("K32").Value = 2344.2075
$ws.Range("M32").Value = -2057.2075

# ARM row 45
$ws.Range("H45").Value = 2177.1
$ws.Range("I45").Value = 1879
$ws.Range("J45").Value = 2624.25
$ws.Range("K45").Value = 1879
$ws.Range("L45").Value = 2624.25
$ws.Range("M45").Value = -1502
$ws.Range("N45").Value = -3378.25

# ARM row 61
$ws.Range("H61").Value = 1829.2222
$ws.Range("I61").Value = 1829.2222
$ws.Range("K61").Value = 1829.2222
$ws.Range("M61").Value = -1617.2222

# ARM row 74
$ws.Range("H74").Value = 1115.9333
$ws.Range("I74").Value = 1117.0714
$ws.Range("J74").Value = 1100
$ws.Range("K74").Value = 1117.0714
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = -243.0714
$ws.Range("N74").Value = -2848

# ARM row 77
$ws.Range("H77").Value = 1115.9333
$ws.Range("I77").Value = 1117.0714
$ws.Range("J77").Value = 1100
$ws.Range("K77").Value = 5585.357
$ws.Range("L77").Value = 5500
$ws.Range("M77").Value = -1217.357
$ws.Range("N77").Value = -14236

# ARM row 122
$ws.Range("H122").Value = 1931.625
$ws.Range("I122").Value = 1850.5
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 5551.5
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -3101.5
$ws.Range("N122").Value = -12398.5

# ARM row 132
$ws.Range("H132").Value = 2244.818
$ws.Range("I132").Value = 2244.818
$ws.Range("K132").Value = 6734.454000000001
$ws.Range("M132").Value = -4204.454000000001

# ARM row 136
$ws.Range("H136").Value = 1829.2222
$ws.Range("I136").Value = 1829.2222
$ws.Range("K136").Value = 5487.6666
$ws.Range("M136").Value = -2937.6666

$ws = $wb.Worksheets.Item("BSM")
# BSM row 82
$ws.Range("H82").Value = 31480.938
$ws.Range("I82").Value = 18707.5
$ws.Range("J82").Value = 39145
$ws.Range("K82").Value = 18707.5
$ws.Range("L82").Value = 39145
$ws.Range("M82").Value = -18324.5
$ws.Range("N82").Value = -39911

# BSM row 85
$ws.Range("H85").Value = 31480.938
$ws.Range("I85").Value = 18707.5
$ws.Range("J85").Value = 39145
$ws.Range("K85").Value = 18707.5
$ws.Range("L85").Value = 39145
$ws.Range("M85").Value = -17381.5
$ws.Range("N85").Value = -41797

# BSM row 86
$ws.Range("H86").Value = 4245.1113
$ws.Range("I86").Value = 3879
$ws.Range("J86").Value = 4820.4287
$ws.Range("K86").Value = 3879
$ws.Range("L86").Value = 4820.4287
$ws.Range("M86").Value = -2756
$ws.Range("N86").Value = -7066.4287

# BSM row 89
$ws.Range("H89").Value = 4245.1113
$ws.Range("I89").Value = 3879
$ws.Range("J89").Value = 4820.4287
$ws.Range("K89").Value = 19395
$ws.Range("L89").Value = 24102.1435
$ws.Range("M89").Value = -13779
$ws.Range("N89").Value = -35334.14350000001

# BSM row 105
$ws.Range("H105").Value = 3280.1
$ws.Range("I105").Value = 3315.875
$ws.Range("J105").Value = 3137
$ws.Range("K105").Value = 3315.875
$ws.Range("L105").Value = 3137
$ws.Range("M105").Value = -1568.875
$ws.Range("N105").Value = -6631

# BSM row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# BSM row 134
$ws.Range("H134").Value = 13692.25
$ws.Range("I134").Value = 13692.25
$ws.Range("K134").Value = 41076.75
$ws.Range("M134").Value = -38541.75

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws.Range("H7").Value = 384667.62
$ws.Range("I7").Value = 52.636364
$ws.Range("K7").Value = 52.636364
$ws.Range("M7").Value = 60.363636

# CRP row 86
$ws.Range("H86").Value = 8712.421
$ws.Range("I86").Value = 8630
$ws.Range("J86").Value = 8772.362999999999
$ws.Range("K86").Value = 8630
$ws.Range("L86").Value = 8772.362999999999
$ws.Range("M86").Value = -7507
$ws.Range("N86").Value = -11018.363

# CRP row 89
$ws.Range("H89").Value = 8712.421
$ws.Range("I89").Value = 8630
$ws.Range("J89").Value = 8772.362999999999
$ws.Range("K89").Value = 43150
$ws.Range("L89").Value = 43861.815
$ws.Range("M89").Value = -37534
$ws.Range("N89").Value = -55093.815

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Range("H5").Value = 1952.2307
$ws.Range("I5").Value = 1954.1111
$ws.Range("J5").Value = 1948
$ws.Range("K5").Value = 5862.3333
$ws.Range("L5").Value = 5844
$ws.Range("M5").Value = -5750.3333
$ws.Range("N5").Value = -6068

# CUL row 29
$ws.Range("H29").Value = 230
$ws.Range("I29").Value = 225
$ws.Range("J29").Value = 231.66667
$ws.Range("K29").Value = 675
$ws.Range("L29").Value = 695.00001
$ws.Range("M29").Value = -398
$ws.Range("N29").Value = -1249.00001

# CUL row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").ClearContents()

# CUL row 38
$ws.Range("H38").Value = 150.2
$ws.Range("I38").Value = 149.625
$ws.Range("K38").Value = 448.875
$ws.Range("M38").Value = -101.875

# CUL row 40
$ws.Range("H40").Value = 597
$ws.Range("J40").Value = 1999.5
$ws.Range("L40").Value = 7998
$ws.Range("N40").Value = -8136

# CUL row 46
$ws.Range("H46").Value = 1043
$ws.Range("I46").Value = 64.5
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 193.5
$ws.Range("L46").Value = 9000
$ws.Range("M46").Value = -102.5
$ws.Range("N46").Value = -9182

# CUL row 69
$ws.Range("H69").Value = 998.5
$ws.Range("J69").Value = 998.5
$ws.Range("L69").Value = 2995.5
$ws.Range("N69").Value = -4617.5

# CUL row 72
$ws.Range("H72").Value = 998.5
$ws.Range("J72").Value = 998.5
$ws.Range("L72").Value = 8986.5
$ws.Range("N72").Value = -17098.5

# CUL row 135
$ws.Range("H135").Value = 1952.2307
$ws.Range("I135").Value = 1954.1111
$ws.Range("J135").Value = 1948
$ws.Range("K135").Value = 17586.9999
$ws.Range("L135").Value = 17532
$ws.Range("M135").Value = -15051.9999
$ws.Range("N135").Value = -22602

$ws = $wb.Worksheets.Item("GSM")
# GSM row 2
$ws.Range("H2").Value = 133.11111
$ws.Range("I2").Value = 106
$ws.Range("K2").Value = 106
$ws.Range("M2").Value = 7

# GSM row 70
$ws.Range("H70").Value = 55562388
$ws.Range("I70").Value = 83339210
$ws.Range("J70").Value = 8749.5
$ws.Range("K70").Value = 83339210
$ws.Range("L70").Value = 8749.5
$ws.Range("M70").Value = -83338940
$ws.Range("N70").Value = -9289.5

# GSM row 73
$ws.Range("H73").Value = 55562388
$ws.Range("I73").Value = 83339210
$ws.Range("J73").Value = 8749.5
$ws.Range("K73").Value = 83339210
$ws.Range("L73").Value = 8749.5
$ws.Range("M73").Value = -83338274
$ws.Range("N73").Value = -10621.5

# GSM row 122
$ws.Range("H122").Value = 2400.9167
$ws.Range("I122").Value = 1353.625
$ws.Range("J122").Value = 4495.5
$ws.Range("K122").Value = 4060.875
$ws.Range("L122").Value = 13486.5
$ws.Range("M122").Value = -1610.875
$ws.Range("N122").Value = -18386.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 22
$ws.Range("H22").Value = 1051.8182
$ws.Range("I22").Value = 957
$ws.Range("K22").Value = 957
$ws.Range("M22").Value = -662

# LTW row 27
$ws.Range("H27").Value = 1051.8182
$ws.Range("I27").Value = 957
$ws.Range("K27").Value = 957
$ws.Range("M27").Value = -850

# LTW row 40
$ws.Range("H40").Value = 2745.25
$ws.Range("I40").Value = 2745.25
$ws.Range("K40").Value = 2745.25
$ws.Range("M40").Value = -2609.25

# LTW row 132
$ws.Range("H132").Value = 2815.0833
$ws.Range("J132").Value = 4000
$ws.Range("L132").Value = 12000
$ws.Range("N132").Value = -17060

# LTW row 136
$ws.Range("H136").Value = 3444.7144
$ws.Range("I136").Value = 3268
$ws.Range("K136").Value = 9804
$ws.Range("M136").Value = -7254

$ws = $wb.Worksheets.Item("WVR")
# WVR row 122
$ws.Range("H122").Value = 2016.3636
$ws.Range("I122").Value = 1748.9445
$ws.Range("K122").Value = 5246.833500000001
$ws.Range("M122").Value = -2796.833500000001

# WVR row 132
$ws.Range("H132").Value = 1414.8334
$ws.Range("I132").Value = 1414.8334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4244.5002
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1714.5002
$ws.Range("N132").ClearContents()

# WVR row 136
$ws.Range("H136").Value = 4095.4722
$ws.Range("I136").Value = 3996.1
$ws.Range("J136").Value = 4592.3335
$ws.Range("K136").Value = 11988.3
$ws.Range("L136").Value = 13777.0005
$ws.Range("M136").Value = -9438.299999999999
$ws.Range("N136").Value = -18877.0005
